# The presentation currently uses the "Integral" theme colors (stored in
# ppt/theme/theme2.xml, the theme actually linked from the slide master /
# presentation). The edit swaps the deck back to the default "Office Theme"
# colors. Apply the Office Theme color scheme to the active theme via the
# PowerPoint object model (Master.Theme.ThemeColorScheme).

function ConvertTo-BGR([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# Office Theme color scheme, in ThemeColorScheme.Item() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-BGR $officeThemeHex[$i - 1]
}
